# Refresh the cryptos price/volume table (GitHub Actions data pull).
# Price cells (column D) are forced back to text via a leading apostrophe and
# then reset to the "Normal" style so the stored number format stays the
# workbook default (no style index) - otherwise Excel would silently parse
# strings like "1.000" or "35.40" as numbers and drop the meaningful
# trailing/internal zeros that the source feed relies on.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.Value = "'30.419.71"
$c.Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  +0.66%  '
$c = $ws.Cells.Item(3, 4)
$c.Value = "'1.869.21"
$c.Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  +0.27%  '
$c = $ws.Cells.Item(5, 4)
$c.Value = "'246.23"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +1.86%  '
$ws.Cells.Item(6, 5).Value = '  +0.03%  '
$c = $ws.Cells.Item(7, 4)
$c.Value = "'0.4739"
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  +0.75%  '
$c = $ws.Cells.Item(8, 4)
$c.Value = "'0.2915"
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  +2.25%  '
$c = $ws.Cells.Item(9, 4)
$c.Value = "'0.06495"
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  +0.39%  '
$ws.Cells.Item(10, 5).Value = '  +6.02%  '
$c = $ws.Cells.Item(11, 4)
$c.Value = "'0.07719"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +0.31%  '
$c = $ws.Cells.Item(12, 4)
$c.Value = "'97.51"
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  +2.70%  '
$c = $ws.Cells.Item(13, 4)
$c.Value = "'0.7384"
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  +8.05%  '
$c = $ws.Cells.Item(14, 4)
$c.Value = "'1.872.16"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +0.40%  '
$c = $ws.Cells.Item(15, 4)
$c.Value = "'5.129"
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +0.84%  '
$c = $ws.Cells.Item(16, 4)
$c.Value = "'273.41"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +1.67%  '
$c = $ws.Cells.Item(17, 4)
$c.Value = "'30.403.52"
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  +0.64%  '
$c = $ws.Cells.Item(18, 4)
$c.Value = "'13.37"
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  +0.05%  '
$c = $ws.Cells.Item(19, 4)
$c.Value = "'0.000007538"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +0.14%  '
$c = $ws.Cells.Item(20, 4)
$c.Value = "'1.000"
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +0.00%  '
$c = $ws.Cells.Item(21, 4)
$c.Value = "'2.117.08"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +0.36%  '
$c = $ws.Cells.Item(22, 4)
$c.Value = "'1.000"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +0.03%  '
$c = $ws.Cells.Item(23, 4)
$c.Value = "'5.223"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +0.44%  '
$c = $ws.Cells.Item(24, 4)
$c.Value = "'6.165"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +0.77%  '
$c = $ws.Cells.Item(25, 4)
$c.Value = "'9.303"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -0.30%  '
$c = $ws.Cells.Item(26, 4)
$c.Value = "'164.23"
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -0.88%  '
$ws.Cells.Item(27, 5).Value = '  -0.04%  '
$ws.Cells.Item(28, 5).Value = '  +1.77%  '
$ws.Cells.Item(29, 2).Value = 'Stellar'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Cells.Item(29, 4)
$c.Value = "'0.09986"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +1.81%  '
$ws.Cells.Item(30, 2).Value = 'Toncoin'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Cells.Item(30, 4)
$c.Value = "'1.367"
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  -0.43%  '
$c = $ws.Cells.Item(31, 4)
$c.Value = "'1.501"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  -0.21%  '
$c = $ws.Cells.Item(32, 4)
$c.Value = "'4.298"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  +1.29%  '
$c = $ws.Cells.Item(33, 4)
$c.Value = "'4.138"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  +3.77%  '
$c = $ws.Cells.Item(34, 4)
$c.Value = "'0.04836"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +2.93%  '
$ws.Cells.Item(35, 5).Value = '  +0.79%  '
$c = $ws.Cells.Item(36, 4)
$c.Value = "'0.6975"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +1.87%  '
$ws.Cells.Item(37, 2).Value = 'Frax'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$c = $ws.Cells.Item(37, 4)
$c.Value = "'0.9995"
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +0.00%  '
$ws.Cells.Item(38, 2).Value = 'HuobiToken'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$c = $ws.Cells.Item(38, 4)
$c.Value = "'2.715"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +0.22%  '
$ws.Cells.Item(39, 2).Value = 'VeChain'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Cells.Item(39, 4)
$c.Value = "'0.01860"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  +0.91%  '
$ws.Cells.Item(40, 2).Value = 'MXToken'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Cells.Item(40, 4)
$c.Value = "'2.743"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +0.07%  '
$ws.Cells.Item(41, 2).Value = 'FraxShare'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Cells.Item(41, 4)
$c.Value = "'6.307"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -0.58%  '
$ws.Cells.Item(42, 2).Value = 'Aave'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Cells.Item(42, 4)
$c.Value = "'73.03"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  +3.70%  '
$ws.Cells.Item(43, 2).Value = 'RenderToken'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Cells.Item(43, 4)
$c.Value = "'1.968"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +4.34%  '
$ws.Cells.Item(44, 2).Value = 'TheSandbox'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$c = $ws.Cells.Item(44, 4)
$c.Value = "'0.4193"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +3.09%  '
$ws.Cells.Item(45, 2).Value = 'PaxDollar'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$c = $ws.Cells.Item(45, 4)
$c.Value = "'1.000"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +0.05%  '
$ws.Cells.Item(46, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Cells.Item(46, 4)
$c.Value = "'0.8337"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  -0.51%  '
$ws.Cells.Item(47, 2).Value = 'Quant'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c = $ws.Cells.Item(47, 4)
$c.Value = "'102.07"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +0.22%  '
$ws.Cells.Item(48, 2).Value = 'EnergySwap'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Cells.Item(48, 4)
$c.Value = "'9.267"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +1.33%  '
$ws.Cells.Item(49, 2).Value = 'Aptos'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c = $ws.Cells.Item(49, 4)
$c.Value = "'6.998"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +0.67%  '
$ws.Cells.Item(50, 2).Value = 'Maker'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$c = $ws.Cells.Item(50, 4)
$c.Value = "'930.21"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +0.85%  '
$ws.Cells.Item(51, 2).Value = 'Elrond'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$c = $ws.Cells.Item(51, 4)
$c.Value = "'35.40"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +2.77%  '
